$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain stored as text even when the value looks numeric,
# then restore the "Normal" style so no stray style index is left on the cell.
function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "66.013.63"
$ws.Range("E2").Value = "  -1.91%  "

Set-TextValue "D3" "3.442.29"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue "D5" "584.15"
$ws.Range("E5").Value = "  -0.70%  "

Set-TextValue "D6" "173.70"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("E7").Value = "  -0.07%  "

Set-TextValue "D8" "0.604"
$ws.Range("E8").Value = "  -0.97%  "

Set-TextValue "D9" "3.441.36"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("E10").Value = "  -3.91%  "

Set-TextValue "D11" "6.93"
$ws.Range("E11").Value = "  -0.35%  "

Set-TextValue "D12" "0.410"
$ws.Range("E12").Value = "  -3.70%  "

Set-TextValue "D13" "4.033.51"
$ws.Range("E13").Value = "  -1.12%  "

$ws.Range("E14").Value = "  +1.28%  "

Set-TextValue "D15" "28.76"
$ws.Range("E15").Value = "  -9.80%  "

Set-TextValue "D16" "66.018.97"
$ws.Range("E16").Value = "  -1.96%  "

Set-TextValue "D17" "0.0000171"
$ws.Range("E17").Value = "  -2.32%  "

Set-TextValue "D18" "3.437.40"
$ws.Range("E18").Value = "  -0.74%  "

Set-TextValue "D19" "5.94"
$ws.Range("E19").Value = "  -2.77%  "

Set-TextValue "D20" "13.86"
$ws.Range("E20").Value = "  -0.54%  "

Set-TextValue "D21" "370.28"
$ws.Range("E21").Value = "  -2.68%  "

Set-TextValue "D22" "7.67"
$ws.Range("E22").Value = "  -2.30%  "

Set-TextValue "D23" "72.43"
$ws.Range("E23").Value = "  +0.90%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("E25").Value = "  +0.27%  "

Set-TextValue "D26" "0.0000122"
$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("E27").Value = "  -2.29%  "

Set-TextValue "D28" "0.177"
$ws.Range("E28").Value = "  +1.34%  "

Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.27%  "

Set-TextValue "D30" "23.66"
$ws.Range("E30").Value = "  -1.52%  "

Set-TextValue "D31" "5.74"
$ws.Range("E31").Value = "  -3.96%  "

$ws.Range("E32").Value = "  -2.54%  "

Set-TextValue "D33" "1.00"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  -5.69%  "

Set-TextValue "D35" "7.01"
$ws.Range("E35").Value = "  -2.87%  "

$ws.Range("E36").Value = "  -1.26%  "

Set-TextValue "D37" "160.74"
$ws.Range("E37").Value = "  +0.17%  "

Set-TextValue "D38" "28.81"
$ws.Range("E38").Value = "  +5.52%  "

Set-TextValue "D39" "0.880"
$ws.Range("E39").Value = "  +0.04%  "

$ws.Range("E40").Value = "  -3.07%  "

Set-TextValue "D41" "2.60"
$ws.Range("E41").Value = "  -1.89%  "

Set-TextValue "D42" "2.769.27"
$ws.Range("E42").Value = "  +2.24%  "

Set-TextValue "D44" "4.46"
$ws.Range("E44").Value = "  -1.02%  "

Set-TextValue "D45" "0.0681"
$ws.Range("E45").Value = "  -2.58%  "

Set-TextValue "D46" "40.23"
$ws.Range("E46").Value = "  -2.38%  "

Set-TextValue "D47" "24.40"
$ws.Range("E47").Value = "  -4.67%  "

$ws.Range("E48").Value = "  -1.78%  "

Set-TextValue "D49" "325.27"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D50" "0.102"
$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "6.26"
$ws.Range("E51").Value = "  +0.05%  "
